# ------------------------------------------------------------------
# q2 complete; updated q2 parameters; q4 results
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ===================== q2 : append new experiment rows =====================
$ws2 = $wb.Worksheets.Item("q2")

# New rows are appended below the existing 8 data rows (rows 2-9), in the
# same chronological order the author typed them, so that new shared
# strings land in the notes column in the same order as the target file.
# A sort (matching the workbook's own sortState) is applied afterwards,
# which re-orders rows 2-25 without touching the shared-string table.
$newRows = @(
    @(100, 100, 0.05, "y", 0.009,  7.06,  3400,  $null),
    @(110, 100, 0.03, "y", 0.009,  8.18,  3700,  $null),
    @(100, 100, 0.03, "y", 0.0097, 8.5,   3700,  $null),
    @(95,  100, 0.05, "y", 0.009,  7.96,  3900,  $null),
    @(105, 100, 0.05, "y", 0.009,  9.31,  4300,  $null),
    @(95,  100, 0.03, "y", 0.009,  10.43, 4900,  $null),
    @(15,  10,  0.05, "n", 1.56,   59,    5800,  $null),
    @(15,  20,  0.05, "n", 0.26,   120,   11700, $null),
    @(25,  20,  0.05, "n", 0.388,  62,    13000, $null),
    @(40,  25,  0.05, "n", 0.38,   435,   27800, "lots of variation even as epochs increase"),
    @(18,  20,  0.03, "n", 0.034,  462,   84800, "was at .05 by 23,000"),
    @(18,  20,  0.05, "n", 0.14,   230,   30000, "stuck at .14"),
    @(14,  20,  0.04, "n", 0.807,  99,    17600, ".808 by 1000"),
    @(16,  25,  0.03, "n", 0.46,   430,   95500, $null),
    @(12,  100, 0.03, "n", 0.072,  150,   143000, $null),
    @(15,  200, 0.05, "n", 0.06,   320,   342000, "lots of variation even as epochs increase")
)

$r = 10
foreach ($row in $newRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -ne $null) {
        $ws2.Cells.Item($r, 8).Value = $row[7]
    }
    $r = $r + 1
}

# Sort the full data range: converge? desc, then time/final loss/epochs asc
# (matches the <sortState> Excel itself records for this range).
$sortRange = $ws2.Range("A2:H26")
$keyD = $ws2.Range("D2:D26")
$keyG = $ws2.Range("G2:G26")
$keyE = $ws2.Range("E2:E26")
$keyF = $ws2.Range("F2:F26")
$ws2.Sort.SortFields.Clear()
$ws2.Sort.SortFields.Add($keyD, [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlDescending) | Out-Null
$ws2.Sort.SortFields.Add($keyG, [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending) | Out-Null
$ws2.Sort.SortFields.Add($keyE, [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending) | Out-Null
$ws2.Sort.SortFields.Add($keyF, [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending) | Out-Null
$ws2.Sort.SetRange($sortRange)
$ws2.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws2.Sort.Apply()

$ws2.Activate()
$ws2.Range("E4").Select()

# ===================== q3 : selection only =====================
$ws3 = $wb.Worksheets.Item("q3")
$ws3.Activate()
$ws3.Range("A1:I1").Select()

# ===================== q4 : new sheet =====================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "q4"

$ws4.Cells.Item(1, 1).Value = "num hidden layers"
$ws4.Cells.Item(1, 2).Value = "hidden layer"
$ws4.Cells.Item(1, 3).Value = "batch"
$ws4.Cells.Item(1, 4).Value = "learning rate"
$ws4.Cells.Item(1, 5).Value = "converge?"
$ws4.Cells.Item(1, 6).Value = "accuracy"
$ws4.Cells.Item(1, 7).Value = "time"
$ws4.Cells.Item(1, 8).Value = "epochs"
$ws4.Cells.Item(1, 9).Value = "notes"

$ws4.Cells.Item(2, 2).Value = 100
$ws4.Cells.Item(2, 3).Value = 25
$ws4.Cells.Item(2, 4).Value = 0.05
$ws4.Cells.Item(2, 5).Value = "y"
$ws4.Cells.Item(2, 6).Value = 85.2
$ws4.Cells.Item(2, 7).Value = 12.69
$ws4.Cells.Item(2, 8).Value = 6

$ws4.Columns.Item(5).ColumnWidth = 9.5

$ws4.Activate()
$ws4.Range("A2").Select()
